# auto:removing some labels from the patient card
#
# Remove the "nick" (Nickname) and "gender_n" (Gender Identity) note rows
# from the patient-card group on the "survey" sheet, and rename the
# settings-sheet form title from "Create a task " to "Task".

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# Delete row 44 first (gender_n) so row 42's (nick) index doesn't move
# before we delete it.
$survey.Rows.Item(44).Delete()
$survey.Rows.Item(42).Delete()

$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "Task"
